$d = $word.ActiveDocument

# 1. Insert a new heading paragraph ("ISOFT Services & Solutions") before the
#    very first paragraph of the document, with the exact pPr/rPr formatting
#    used by the target document (teal color, 14pt/28 half-points, and the
#    5770/710 twip indent).
$insertionPoint = $d.Range(0, 0)
$newParaXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:spacing w:after="290" w:line="259" w:lineRule="auto"/>
              <w:ind w:left="5770" w:firstLine="710"/>
              <w:rPr>
                <w:color w:val="2A7B88"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:color w:val="2A7B88"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
              </w:rPr>
              <w:t>ISOFT Services &amp; Solutions</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
[void]$insertionPoint.InsertXML($newParaXml)

# 2. The paragraph that used to be first ("Nagulu Syed ...") now sits second;
#    give it an explicit zero indent so it matches the diff.
$origFirstPara = $d.Paragraphs.Item(2)
$origFirstPara.Format.LeftIndent = 0
$origFirstPara.Format.FirstLineIndent = 0

# 3. Shrink the document margins from 1440/1080 to a uniform 720 twips
#    (36pt) on every section.
foreach ($sec in $d.Sections) {
    $sec.PageSetup.TopMargin = 36
    $sec.PageSetup.BottomMargin = 36
    $sec.PageSetup.LeftMargin = 36
    $sec.PageSetup.RightMargin = 36
}

Write-Output "done"
